# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.265.30'
$ws.Range("E2").Value = '  +4.98%  '
$ws.Range("D3").Value = '3.776.63'
$ws.Range("E3").Value = '  +22.28%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'618.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.19%  '
$ws.Range("D6").Value = "'177.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.14%  '
$ws.Range("D7").Value = '3.772.68'
$ws.Range("E7").Value = '  +22.19%  '
$ws.Range("D9").Value = "'0.551"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.93%  '
$ws.Range("E10").Value = '  +12.12%  '
$ws.Range("D11").Value = "'6.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.93%  '
$ws.Range("D12").Value = "'0.506"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.84%  '
$ws.Range("D13").Value = "'40.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +13.52%  '
$ws.Range("D14").Value = "'0.0000259"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.13%  '
$ws.Range("D15").Value = '4.410.69'
$ws.Range("E15").Value = '  +22.39%  '
$ws.Range("D16").Value = '3.783.24'
$ws.Range("E16").Value = '  +22.56%  '
$ws.Range("D17").Value = '70.437.20'
$ws.Range("E17").Value = '  +5.30%  '
$ws.Range("E18").Value = '  +1.47%  '
$ws.Range("D19").Value = "'7.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.10%  '
$ws.Range("D20").Value = "'525.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.78%  '
$ws.Range("D21").Value = "'16.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.08%  '
$ws.Range("D22").Value = "'9.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +24.37%  '
$ws.Range("D23").Value = "'0.748"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.29%  '
$ws.Range("D24").Value = "'88.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.79%  '
$ws.Range("E25").Value = '  +11.46%  '
$ws.Range("D26").Value = "'13.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.32%  '
$ws.Range("D27").Value = "'11.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.67%  '
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").Value = "'0.0000122"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +31.16%  '
$ws.Range("D30").Value = "'2.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.22%  '
$ws.Range("D31").Value = "'2.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.26%  '
$ws.Range("D32").Value = "'8.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.40%  '
$ws.Range("D33").Value = "'32.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +15.75%  '
$ws.Range("E34").Value = '  +4.02%  '
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").Value = "'6.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +11.69%  '
$ws.Range("D37").Value = "'1.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.27%  '
$ws.Range("D38").Value = "'0.344"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.95%  '
$ws.Range("D39").Value = "'2.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.98%  '
$ws.Range("D40").Value = "'0.134"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.94%  '
$ws.Range("D41").Value = "'51.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.28%  '
$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D42").Value = "'8.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.42%  '
$ws.Range("B43").Value = 'Arweave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D43").Value = "'44.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.71%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = "'429.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +16.80%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '3.141.36'
$ws.Range("E45").Value = '  +13.27%  '
$ws.Range("E46").Value = '  +4.33%  '
$ws.Range("E47").Value = '  +8.16%  '
$ws.Range("D48").Value = "'27.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.08%  '
$ws.Range("D49").Value = "'139.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.06%  '
$ws.Range("D50").Value = "'2.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +11.53%  '
$ws.Range("E51").Value = '  +0.01%  '
